{"js": "// Locate the paragraph that announces the loan amount, e.g.:\n//   \"V\u011b\u0159itel t\u00edmto poskytuje dlu\u017en\u00edkovi pen\u011b\u017eitou p\u016fj\u010dku ve v\u00fd\u0161i:\"\n//   <line break>\n//   \"dvacet tis\u00edc korun \u010desk\u00fdch (20 000 K\u010d).\"\n// and rewrite it so the amount in parentheses becomes the\n// \"[[AMOUNT_1]]\" placeholder, collapsing the paragraph's several runs\n// (plain text / break / bold amount / trailing period) into a single,\n// non-bold run that still contains the line break.\nconst body = context.document.body;\n\nconst results = body.search(\n  \"V\u011b\u0159itel t\u00edmto poskytuje dlu\u017en\u00edkovi pen\u011b\u017eitou p\u016fj\u010dku ve v\u00fd\u0161i\",\n  { matchCase: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst paragraph = results.items[0].paragraphs.getFirst();\nparagraph.load(\"text\");\nawait context.sync();\n\n// Replace \"(20 000 K\u010d)\" (or whatever amount is currently present inside\n// the parentheses) with the \"[[AMOUNT_1]]\" placeholder while keeping the\n// rest of the paragraph text \u2014 including the embedded line break \u2014\n// untouched.\nconst updatedText = paragraph.text.replace(/\\([^()]*\\)(?=\\.?$)/, \"([[AMOUNT_1]])\");\n\nparagraph.clear();\nparagraph.insertText(updatedText, \"Start\");\nawait context.sync();\n", "ps1": "# Locate the paragraph that announces the loan amount, e.g.:\n#   \"Veritel timto poskytuje dluznikovi penezitou pujcku ve vysi:\"\n#   <line break>\n#   \"dvacet tisic korun ceskych (20 000 Kc).\"\n# and rewrite it so the amount in parentheses becomes the\n# \"[[AMOUNT_1]]\" placeholder, collapsing the paragraph's several runs\n# (plain text / break / bold amount / trailing period) into a single,\n# non-bold run that still contains the line break.\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"V\u011b\u0159itel t\u00edmto poskytuje dlu\u017en\u00edkovi pen\u011b\u017eitou p\u016fj\u010dku ve v\u00fd\u0161i\")\nif (-not $found) {\n    throw \"Target paragraph not found\"\n}\n\n$paragraph = $searchRange.Paragraphs(1)\n$pStart = $paragraph.Range.Start\n$pEnd = $paragraph.Range.End\n\n# Range covering the paragraph's content only (excludes the trailing\n# paragraph mark), so assigning .Text replaces every run in the\n# paragraph with a single new plain run instead of just overwriting\n# part of the first run.\n$target = $d.Range($pStart, $pEnd - 1)\n\n# Swap out whatever amount currently sits inside the parentheses\n# (e.g. \"(20 000 K\u010d)\") for the \"[[AMOUNT_1]]\" placeholder, keeping the\n# rest of the paragraph text -- including the embedded line break --\n# untouched.\n$target.Text = $target.Text -replace '\\([^()]*\\)(?=\\.?$)', '([[AMOUNT_1]])'\n"}
